$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 128, shifting rows 128:249 down to 129:250.
# Copy row 128 down first so the new row 128 inherits the same formatting/
# repeated columns (A,B,C,E,F,G,H,N,O,Q,R) as the row being displaced.
$ws.Rows("128").Copy()
$ws.Rows("128").Insert()

# Now set the new row 128 values that differ from the old row 128
$ws.Range("D128").Value = 44658
$ws.Range("I128").Value = "Primera"
$ws.Range("J128").Value = 2500
$ws.Range("K128").Value = 800
$ws.Range("L128").Value = 800
$ws.Range("M128").Value = 800
$ws.Range("P128").Value = 800
